# Re-attribute tracked changes made in this session to "vchoi" so the
# revision marks line up with the document's existing reviewer.
$word.UserName = "vchoi"

$d = $word.ActiveDocument

# Track Changes should already be on for this document, but make sure -
# the inserted text must show up as a w:ins revision, not a plain edit.
$d.TrackRevisions = $true

# Find the anchor text "The next section " and collapse the found range
# to its end, so we insert right after it (before "provides a").
$anchor = $d.Content
$found = $anchor.Find.Execute("The next section ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text 'The next section '"
}
$anchor.Collapse(0)

# Insert "User Guide " (tracked as an insertion) right after the anchor.
$newText = $anchor.Duplicate
$newText.InsertAfter("User Guide ")

# Italicize just the words "User Guide" (not the trailing space) -
# this records an rPrChange on that run capturing the prior (non-italic)
# formatting, matching how the author highlighted the phrase after typing it.
$phrase = $d.Range($newText.Start, $newText.Start + 10)
$phrase.Italic = 1
